# Bugfix sur la masse des LEDs. Recompile OK. Ajout des outputjobs demandes
#
# Applies the BOM updates: new footprints for several parts, LED/Diode row
# swap (rows 6 & 7), resistor row reorder + renumbering (rows 13 & 14),
# potentiometer value/footprint/ref updates (row 12), and the saved window
# size in the workbook view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (10uF capacitor): footprint RAD-0.3 -> CAPC2012N
$ws.Range("G4").Value = "CAPC2012N"

# --- Rows 6 & 7 swap content: "GF1A" diode <-> "Green LED" LED -----------
# Row 6 becomes the Green LED line (previously row 7)
$ws.Range("A6").Value = "Green LED"
$ws.Range("B6").Value = "Typical RED, GREEN, YELLOW, AMBER GaAs LED"
$ws.Range("C6").Value = "LED2"
$ws.Range("D6").Value = 5
$ws.Range("F6").Value = "KINGBRIGHT"
$ws.Range("G6").Value = "LED 0805"
$ws.Range("I6").Value = "KPT-2012SGC"
$ws.Range("J6").Value = "2099239"
$ws.Range("K6").Value = "D1_LED0, D1_LED1, D1_LED2, D1_LED3, D4"

# Row 7 becomes the GF1A diode line (previously row 6)
$ws.Range("A7").Value = "GF1A"
$ws.Range("B7").Value = "Default Diode"
$ws.Range("C7").Value = "Diode"
$ws.Range("D7").Value = 2
$ws.Range("F7").Value = "VISHAY"
$ws.Range("G7").Value = "SMA/DO-214AC_21"
$ws.Range("I7").Value = "GF1A-E3/67A"
$ws.Range("J7").Value = "9549560"
$ws.Range("K7").Value = "D2, D3"

# --- Row 8 (Inductor FB1): footprint INDC1608L -> INDC1608AN --------------
$ws.Range("G8").Value = "INDC1608AN"

# --- Row 9 (Program / Debug J1): footprint HDR1X6H -> RJ12_90 -------------
$ws.Range("G9").Value = "RJ12_90"

# --- Rows 10 & 11 (Header 2 P1/P2): footprint HDR1X2 -> PHOENIX_1985195 ---
$ws.Range("G10").Value = "PHOENIX_1985195"
$ws.Range("G11").Value = "PHOENIX_1985195"

# --- Row 12 (Potentiometer R1): value/fabricant/footprint/refs updated ----
$ws.Range("A12").Value = "10k"
$ws.Range("E12").Value = "10K"
$ws.Range("F12").Value = "Bourns"
$ws.Range("G12").Value = "PDB181-K415K-102A2"
$ws.Range("I12").Value = "''PDB181-K415K-102A2"
$ws.Range("J12").Value = "''1823540"

# --- Rows 13 & 14 (Resistor) reorder & renumber ---------------------------
# Row 13 becomes the 330r / R2_LED... line (previously row 14)
$ws.Range("A13").Value = "330r"
$ws.Range("D13").Value = 4
$ws.Range("G13").Value = "RESC1608N"
$ws.Range("K13").Value = "R2_LED0, R2_LED1, R2_LED2, R2_LED3"

# Row 14 becomes the 4.7k / R3 line (previously row 13)
$ws.Range("A14").Value = "4.7k"
$ws.Range("D14").Value = 1
$ws.Range("G14").Value = "RESC1608N"
$ws.Range("K14").Value = "R3"

# --- Row 15 (Resistor R4): footprint AXIAL-0.3 -> RESC1608N ---------------
$ws.Range("G15").Value = "RESC1608N"

# --- Row 16 (PIC32MX U1): footprint blank -> QFP50P1200X1200X120-64 -------
$ws.Range("G16").Value = "QFP50P1200X1200X120-64"

# --- Row 17 (AP1117 U2): footprint D2PAK_M -> TD03B_N ---------------------
$ws.Range("G17").Value = "TD03B_N"

# --- Row 18 (Crystal Oscillator Y1): footprint R38 -> HC49/4H SMX CRYSTAL -
$ws.Range("G18").Value = "HC49/4H SMX CRYSTAL"

# --- Saved window size in the workbook view -------------------------------
$excel.WindowState = -4143
$win = $excel.Windows.Item(1)
$win.Width = 21390
$win.Height = 12405
